# Actualización de horarios Línea 141 - 876
# Actualiza las 3 hojas (LP1912, LP1912-215, 6203-6173) con el nuevo scrap
# de las 04:13:31, insertando las nuevas llegadas y recalculando los
# "Minutos" / "Total filas" / "Última actualización" de cada hoja.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:13:31"
$ws1.Range("A3").Value = "Total filas: 15"

$sheet1Rows = @(
    @("03:52:04", "04:01", "81_EL PELIGRO",  9,   "LP1912"),
    @("04:13:31", "04:46", "215_EL PELIGRO", 33,  "LP1912"),
    @("03:52:04", "04:46", "215A_EL PATO",   54,  "LP1912"),
    @("04:13:31", "04:53", "11_ETCHEVERRY",  40,  "LP1912"),
    @("04:13:31", "05:11", "17_ROMERO",      58,  "LP1912"),
    @("03:52:04", "05:16", "17_ROMERO",      84,  "LP1912"),
    @("04:13:31", "05:22", "23_HERNANDEZ",   69,  "LP1912"),
    @("04:13:31", "05:31", "81_EL PELIGRO",  78,  "LP1912"),
    @("03:52:04", "05:35", "215B_EL PATO",   103, "LP1912"),
    @("03:52:04", "05:46", "15_ABASTO",      114, "LP1912"),
    @("04:13:31", "05:50", "14_ABASTO",      97,  "LP1912"),
    @("04:13:31", "05:52", "17_ROMERO",      99,  "LP1912"),
    @("04:13:31", "06:01", "16_SANTA ANA",   108, "LP1912"),
    @("04:13:31", "06:03", "10_OLMOS",       110, "LP1912"),
    @("04:13:31", "06:11", "215A_EL PATO",   118, "LP1912")
)

$r = 6
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:13:31"
$ws2.Range("A3").Value = "Total filas: 4"

$sheet2Rows = @(
    @("04:13:31", "04:46", "215_EL PELIGRO", 33,  "LP1912"),
    @("03:52:04", "04:46", "215A_EL PATO",   54,  "LP1912"),
    @("03:52:04", "05:35", "215B_EL PATO",   103, "LP1912"),
    @("04:13:31", "06:11", "215A_EL PATO",   118, "LP1912")
)

$r = 6
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Hoja 3: 6203-6173 (solo se actualiza la hora de "Última actualización")
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:13:31"
